$wb = $excel.ActiveWorkbook

# --- Fix the instruction text on "FlightsPosDefaultDate" so it references
#     its own tag name (FlightsDefaultDtPos) instead of the stale one ---
$wsDefaultDate = $wb.Worksheets.Item("FlightsPosDefaultDate")
$wsDefaultDate.Range("I9").Value = "Instruction - Put the start tag (FlightsDefaultDtPos) before the first column and the first row to be used and the end tag (FlightsDefaultDtPos) after the last column and last row to be used.`n`nOnly use positive scenarios here."

# --- Rename the "OneWayFlightsPosDate" sheet to "FlightsPosDate" ---
$wsPosDate = $wb.Worksheets.Item("OneWayFlightsPosDate")
$wsPosDate.Name = "FlightsPosDate"

# --- Update the tag / instruction text on the renamed sheet ---
# The start/end tag cells (A1 and D7) change from "OneWayFlightsPosDate" to
# "FlightsPosWithDate" - this is the new return-flight scenario.
$wsPosDate.Range("A1").Value = "FlightsPosWithDate"
$wsPosDate.Range("D7").Value = "FlightsPosWithDate"

# The instruction cell (I9) is updated to reference the new tag name.
$wsPosDate.Range("I9").Value = "Instruction - Put the start tag (FlightsPosWithDate) before the first column and the first row to be used and the end tag (FlightsPosWithDate) after the last column and last row to be used.`n`nOnly use positive scenarios here."

# --- Move the active tab / selection from "OneWayFlightsPosAllOptions" to
#     the renamed "FlightsPosDate" sheet, with a new selected cell ---
$wsPosDate.Activate()
$wsPosDate.Range("E9").Select()
